# This script updates the cryptocurrency price/volume table in Sheet1
# to reflect refreshed data (GitHub Actions scheduled update).
# For row 50/51 the two coins (PEPE / BitcoinSV) also swap places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.816.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.29%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.409.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.15%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "412.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.62%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.724"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.19%  "

# Row 10
$ws.Range("E10").Value = "  -5.52%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "

# Row 12
$ws.Range("E12").Value = "  -2.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.72%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.948.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.32%  "

# Row 15
$ws.Range("E15").Value = "  -0.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.68%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.420.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.74%  "

# Row 19
$ws.Range("E19").Value = "  +0.68%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "61.845.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.21%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "479.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +17.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.52%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.85%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.56%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.30%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.96%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.12%  "

# Row 29
$ws.Range("E29").Value = "  +2.66%  "

# Row 30
$ws.Range("E30").Value = "  +0.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.41%  "

# Row 32
$ws.Range("E32").Value = "  -1.86%  "

# Row 33
$ws.Range("E33").Value = "  -3.46%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.13%  "

# Row 35
$ws.Range("E35").Value = "  -0.66%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.93%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0486"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.52%  "

# Row 38
$ws.Range("E38").Value = "  -0.01%  "

# Row 39
$ws.Range("E39").Value = "  +2.93%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.322"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.48%  "

# Row 42
$ws.Range("E42").Value = "  -0.20%  "

# Row 43
$ws.Range("E43").Value = "  -0.55%  "

# Row 44
$ws.Range("E44").Value = "  +4.42%  "

# Row 45
$ws.Range("E45").Value = "  +6.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.08%  "

# Row 47
$ws.Range("E47").Value = "  +19.39%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.78%  "

# Row 50
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0519"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.79%  "

# Row 51
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.13%  "
